$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B92: was stored as text "2", should become a real number 2 ---
$ws.Range("B92").Value = 2

# --- New row 93 ---
$ws.Range("A93").Value = "Ruilin"

# B93 must stay a text string "4" (not get auto-converted to a number),
# so force the cell to text format before writing the value, then strip
# the formatting change back off so no stray style sticks to the cell.
$ws.Range("B93").NumberFormat = "@"
$ws.Range("B93").Value = "4"
$ws.Range("B93").ClearFormats()

$ws.Range("C93").Value = "thank"
$ws.Range("D93").Value = "ACK"
$ws.Range("E93").Value = "OTH"
$ws.Range("F93").Value = "c39fead7-b272-4988-9907-50ea12305918"
$ws.Range("G93").Value = "HknbyQbC-_annotated.xlsx"
$ws.Range("H93").Value = "We thank the reviewer for the thoughtful comments and suggestions."

Write-Output "applied annotation update for Ruilin row 93"
